# Region II_ELECTRIFICATION.xlsx edit
# - Insert two new columns of data (AS, AT) between the existing "Unnamed: 43" (AR)
#   column and the "Status as of July 4, 2025" (old AS, now AU) column.
# - For most data rows, the old AR value (a status string) moves one column to
#   the right, landing in AT (column 46).
# - For a handful of rows ("BBM ..." remarks), the AR cell is overwritten with a
#   new remark value instead of being cleared, and AT still receives the old
#   status value.
# - For the first block of rows (2-21), a new AS value "ongrid" is written
#   and the old status moves to AT.
# - Update the dimension (handled automatically by Excel once all the cells
#   are written) and the dropdown data validation range (AS2:AS378 -> AU2:AU378).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Header row (row 1)
#    AR1 ("Unnamed: 43") keeps its text/style.
#    New AS1 = "Unnamed: 44" (bold/bordered header style, copied from AR1)
#    New AT1 = "Unnamed: 45" (same header style)
#    Old AS1 ("Status as of July 4, 2025") becomes AU1 (no special style)
# ---------------------------------------------------------------------------
$oldAS1 = $ws.Cells.Item(1, 45).Value()   # "Status as of July 4, 2025"

$ws.Cells.Item(1, 45).Value = "Unnamed: 44"   # AS1
$ws.Cells.Item(1, 46).Value = "Unnamed: 45"   # AT1
$ws.Cells.Item(1, 47).Value = $oldAS1         # AU1

# Copy the header formatting (bold font, border, centered/top alignment) from
# AR1 onto the two newly-inserted header cells AS1:AT1.
$ws.Range("AR1").Copy()
$ws.Range("AS1:AT1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Data rows (2-378)
# ---------------------------------------------------------------------------

# Rows 2-21: new AS value is "ongrid"; old AR value moves to AT; AR cleared.
$ongridRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21)

# Rows that get a new remark written into AR (overwriting the old status,
# which moves to AT instead of being lost).
$bbmRemarks = @{
    33  = "BBM 2024 SOLAR"
    36  = "BBM 2024 SOLAR"
    37  = "BBM 2024 SOLAR"
    42  = "BBM 2024 SOLAR"
    52  = "BBM 2024 SOLAR"
    53  = "BBM 2022"
    54  = "BBM 2024 SOLAR"
    56  = "BBM 2024 SOLAR"
    63  = "BBM 2024 SOLAR"
    64  = "BBM 2024 SOLAR"
    66  = "BBM 2024 SOLAR"
    100 = "BBM 2022"
    181 = "BBM 2024 SOLAR"
    182 = "BBM 2024 SOLAR"
    197 = "BBM 2022"
    340 = "BBM 2023 UPGRADE"
    352 = "BBM 2025 UPGRADE"
    358 = "BBM 2023 UPGRADE"
    359 = "BBM 2023 UPGRADE"
    360 = "BBM 2024 UPGRADE"
    361 = "bbm 2023 SOLAR"
    362 = "bbm 2023 SOLAR"
    364 = "BBM 2025 ONGRID"
    365 = "BBM 2025 ONGRID"
    368 = "BBM 2025 ONGRID"
    369 = "BBM 2024 UPGRADE"
    370 = "BBM 2025 SOLAR"
    371 = "BBM 2024 UPGRADE"
    372 = "BBM 2025 ONGRID"
    373 = "BBM 2024 UPGRADE"
    374 = "BBM 2024 UPGRADE"
    375 = "BBM 2025 SOLAR"
}

for ($r = 2; $r -le 378; $r++) {
    $oldAR = $ws.Cells.Item($r, 44).Value()   # AR<r>

    if ($ongridRows -contains $r) {
        $ws.Cells.Item($r, 46).Value = $oldAR      # AT<r> = old status
        $ws.Cells.Item($r, 45).Value = "ongrid"     # AS<r> = "ongrid"
        $ws.Cells.Item($r, 44).ClearContents()      # clear AR<r>
    }
    elseif ($bbmRemarks.ContainsKey($r)) {
        $ws.Cells.Item($r, 46).Value = $oldAR              # AT<r> = old status
        $ws.Cells.Item($r, 44).Value = $bbmRemarks[$r]      # AR<r> = remark
    }
    else {
        $ws.Cells.Item($r, 46).Value = $oldAR   # AT<r> = old status
        $ws.Cells.Item($r, 44).ClearContents()  # clear AR<r>
    }
}

# ---------------------------------------------------------------------------
# 3. Data validation: move the dropdown list from AS2:AS378 to AU2:AU378
# ---------------------------------------------------------------------------
$ws.Range("AS2:AS378").Validation.Delete()
$ws.Range("AU2:AU378").Validation.Add(3, 1, 7, "=DropdownOptions!`$A`$1:`$A`$7")
$newValidation = $ws.Range("AU2").Validation
$newValidation.IgnoreBlank = $true
$newValidation.ShowInput = $false
$newValidation.ShowError = $false
